# Daily attendance processing - 2026-01-10 10:35:01
# Re-orders the "Recorded By" (column G) entries so that "System" is
# listed immediately after any literal lowercase "system" token (or first,
# when there isn't one) for rows that were recorded together with
# dnasr281@gmail.com or backup@backdoor.com.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Text

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "backup@backdoor.com, System") {
        $cell.Value = "System, backup@backdoor.com"
    }
    elseif ($val -eq "system, backup@backdoor.com, System") {
        $cell.Value = "system, System, backup@backdoor.com"
    }
}
